# Update "Wellness" tracker: append 15 new daily-report rows (586-600) for
# the 2025-12-0x session (serial date 45980), mirroring the existing data
# rows, plus extend the I-column "Charge" formula down through row 600.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Stamp out new rows by cloning the formatting of an existing row ---
# Row 580 has a non-empty "Localisation douleur" (G) cell -> use it as the
# template for new rows that also have a G value.
# Row 582 has an empty G cell (style s="2") -> use it as the template for
# new rows whose G value is blank.
$templateWithG    = $ws.Range("A580:I580")
$templateNoG      = $ws.Range("A582:I582")

$rowsWithG  = @(587,588,590,591,592,593,594,595,597,599,600)
$rowsWithoutG = @(586,589,596,598)

foreach ($r in $rowsWithG) {
    $templateWithG.Copy($ws.Range("A$r`:I$r"))
}
foreach ($r in $rowsWithoutG) {
    $templateNoG.Copy($ws.Range("A$r`:I$r"))
}

# --- 2. Fill in the actual values for each new row ---
# Columns: A=Date, B=Nom du joueur, C=Volume, D=Intensité, E=Fatigue,
#          F=Douleur, G=Localisation douleur, H=Plaisir, I=Charge (=C*D)
$data = @(
    @{R=586; B="Amir Etien";        C=70; D=6; E=6; F=0; G="";              H=3},
    @{R=587; B="Ilyes Boughanmi";   C=70; D=5; E=2; F=2; G="Genou";         H=10},
    @{R=588; B="Omar Benyounes";    C=70; D=6; E=7; F=4; G="Ichios ";       H=5},
    @{R=589; B="Karim Belmahi";     C=70; D=6; E=7; F=0; G="";              H=10},
    @{R=590; B="Yoan Zouma";        C=70; D=4; E=7; F=7; G="Ischio";        H=5},
    @{R=591; B="Kamal Bafounta";    C=70; D=7; E=5; F=1; G="Genou";         H=6},
    @{R=592; B="Yoann Martelat";    C=70; D=4; E=6; F=2; G="Genou";         H=5},
    @{R=593; B="Jeremie Laurent";   C=70; D=8; E=5; F=1; G="Courbatures";   H=6},
    @{R=594; B="Levy Ndoutoume";    C=70; D=6; E=6; F=3; G="Ischio";        H=2},
    @{R=595; B="Emmanuel Valey";    C=70; D=6; E=6; F=3; G="Ischio";        H=8},
    @{R=596; B="Ilan Ihaddadene";   C=70; D=7; E=5; F=0; G="";              H=10},
    @{R=597; B="Karahali Souaré";   C=70; D=6; E=6; F=7; G="Cheville";      H=0},
    @{R=598; B="Mattheo Haon";      C=70; D=7; E=7; F=0; G="";              H=8},
    @{R=599; B="Sofiane Belle";     C=70; D=7; E=8; F=8; G="Partout";       H=3},
    @{R=600; B="Hedi Nasri";        C=70; D=7; E=7; F=3; G="Hanche";        H=7}
)

foreach ($row in $data) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = 45980
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    if ($row.G -ne "") {
        $ws.Cells.Item($r, 7).Value = $row.G
    }
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Formula = "=C$r*D$r"
}

# --- 3. Update the sheet view state to match the grown range ---
# (topLeftCell scroll-position isn't exposed on this COM surface; the
# active-cell selection is, so at least that part of the view state
# tracks the real edit.)
$ws.Range("A571").Select() | Out-Null
$ws.Range("K595").Select() | Out-Null
